# Daily attendance processing - 2025-12-06 13:34:15
# Normalizes the "Recorded By" column (G) so that "System" is always listed
# first among the recorders, preserving the relative order of the other
# recorder names/emails.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $raw = $cell.Value2

    if ($raw -eq $null) { continue }

    $parts = $raw -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    # Case-sensitive lookup for an exact "System" entry. Note: the -eq/-ceq
    # operators in this environment are case-insensitive, so .Equals() (a
    # true .NET, case-sensitive comparison) is used instead.
    $hasExactSystem = $false
    foreach ($p in $trimmed) {
        if ($p.Equals("System")) { $hasExactSystem = $true }
    }
    if (-not $hasExactSystem) { continue }

    $rest = @()
    $removedOne = $false
    foreach ($p in $trimmed) {
        if ((-not $removedOne) -and $p.Equals("System")) {
            $removedOne = $true
        } else {
            $rest += $p
        }
    }

    $newParts = @("System") + $rest
    $newValue = [string]::Join(", ", $newParts)

    if (-not $newValue.Equals($raw)) {
        $cell.Value = $newValue
    }
}
